# Update the "Update automatically" date/time placeholder shown on the
# slide master and every slide layout (cached text for the
# datetimeFigureOut field) from 28/4/2016 -> 3/5/2016, and bump the
# version number shown in the version-history diagram on slide 1 from
# 8.6.4 -> 8.7.1.

$p = $ppt.ActivePresentation

$oldDate = "28/4/2016"
$newDate = "3/5/2016"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own Date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout ("custom layout") has its own cached copy of the
# date placeholder text.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Version-history diagram on slide 1: bump the version label.
$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shp = $slide.Shapes.Item($si)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "8.6.4") {
            $shp.TextFrame.TextRange.Text = "8.7.1"
        }
    }
}
